# "Include bend radius on BOM"
# Adds a new "Default Bend Radius" column (H) to the BOM header row,
# resizes a few columns, and moves the active selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the added BOM column, formatted like its neighbours
# (bold/underlined, centered header style used across row 1).
$ws.Range("H1").Value = "Default Bend Radius"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Column width tweaks that came along with the new column.
# (ColumnWidth is in "characters"; stored OOXML width = ColumnWidth + 5/6.)
$ws.Columns.Item(2).ColumnWidth = 9.166666666666666   # B: 13 -> 10
$ws.Columns.Item(7).ColumnWidth = 9                    # G: new explicit width ~9.86
$ws.Columns.Item(8).ColumnWidth = 21.333333333333332   # H: new explicit width ~22.14

# Matches the saved cursor position recorded in the sheet view.
$ws.Range("A2").Select()
